# Update cryptos list with latest price/volume snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.354.99"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "2.935.26"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'593.82"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "'145.18"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.505"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  +2.85%  "
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "'33.68"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "3.420.18"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "61.293.86"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "'6.75"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "2.935.17"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'434.02"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").Value = "'13.52"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "'0.682"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "'7.15"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "'82.06"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").Value = "'2.22"
$ws.Range("D26").Value = "'11.89"
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "'2.24"
$ws.Range("E28").Value = "  -4.01%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'7.02"
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("D32").Value = "'26.78"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").Value = "0.0₃0891"
$ws.Range("E34").Value = "  +3.24%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").Value = "'3.02"
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("D38").Value = "'2.02"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "'8.65"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("D41").Value = "'42.60"
$ws.Range("E41").Value = "  +7.09%  "
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "'371.84"
$ws.Range("E44").Value = "  -3.04%  "
$ws.Range("D45").Value = "2.706.30"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "'133.66"
$ws.Range("E46").Value = "  +3.26%  "
$ws.Range("D48").Value = "'23.99"
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("D49").Value = "'0.106"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("E51").Value = "  -0.13%  "
